# The source data contains a weekly series of Albahaca (basil) price
# observations for "Vega Central Mapocho de Santiago". A new weekly
# observation is inserted right before the existing row that currently
# sits at row 628, pushing all subsequent rows down by one.
#
# This mirrors inserting a new record into the dataset (commit message:
# "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 628; Excel shifts rows 628:661 down to
# 629:662 and copies formatting (incl. the date style) from the row above.
$ws.Rows.Item(628).Insert()

# Populate the newly inserted row 628 with the new weekly observation.
$ws.Range("A628").Value2 = 9
$ws.Range("B628").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C628").Value2 = "Metropolitana"
$ws.Range("D628").Value2 = 45267
$ws.Range("E628").Value2 = 13
$ws.Range("F628").Value2 = 100112052
$ws.Range("G628").Value2 = "Albahaca"
$ws.Range("H628").Value2 = "Sin especificar"
$ws.Range("I628").Value2 = "Primera"
$ws.Range("J628").Value2 = 160
$ws.Range("K628").Value2 = 5000
$ws.Range("L628").Value2 = 6000
$ws.Range("M628").Value2 = 5500
$ws.Range("N628").Value2 = "`$/docena de matas"
$ws.Range("O628").Value2 = "Provincia de Chacabuco"
$ws.Range("P628").Value2 = 917
$ws.Range("Q628").Value2 = 6
$ws.Range("R628").Value2 = "Hortaliza"
